$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in columns AD, AE, AF of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting (bold, centered, thin border) used by the
# existing header row by copying the style from an adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the team record (Wins/Losses/Ties) for every player row (2-41).
$ws.Range("AD2:AD41").Value = 80
$ws.Range("AE2:AE41").Value = 81
$ws.Range("AF2:AF41").Value = 0
